# #5: property boat&car done
# Extend the "汽車" (car) sheet with the same property/legislator metadata
# columns already present on the 土地 (land) and 建物 (building) sheets, and
# add a "capacity" column header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Header row (row 1) ----
# New header label in column C: "capacity"
$ws.Cells.Item(1, 3).Value = "capacity"

# New trailing header labels in columns H..N
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Match the bold / bordered / centered look of the existing header cells
$headerNew = $ws.Range("H1:N1")
$headerNew.Font.Bold = $true
$headerNew.Borders.LineStyle = 1
$headerNew.HorizontalAlignment = -4108
$headerNew.VerticalAlignment = -4160

# ---- Data rows 2-4 ----
$indexValues = @(36, 37, 38)

for ($i = 0; $i -lt $indexValues.Length; $i++) {
    $r = $i + 2

    $ws.Cells.Item($r, 8).Value = "land"
    $ws.Cells.Item($r, 9).Value = "normal"

    # Force text so the engine doesn't reinterpret the ISO-like date as a
    # date serial number.
    $ws.Cells.Item($r, 10).NumberFormat = "@"
    $ws.Cells.Item($r, 10).Value = "2012-04-12"

    $ws.Cells.Item($r, 11).Value = "盧嘉辰"
    $ws.Cells.Item($r, 12).Value = 1715
    $ws.Cells.Item($r, 13).Value = "tmp79201"
    $ws.Cells.Item($r, 14).Value = $indexValues[$i]
}
